$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.567.43'
$ws.Range('E2').Value = '  -0.96%  '
$ws.Range('D3').Value = '2.060.08'
$ws.Range('E3').Value = '  +1.07%  '
$ws.Range('E4').Value = '  +0.06%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '243.10'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.66%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.666'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +1.37%  '
$ws.Range('E7').Value = '  +0.03%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '54.46'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -5.48%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '58.59'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -0.81%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.361'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -3.57%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0751'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -3.04%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.938'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +6.79%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '14.72'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -3.78%  '
$ws.Range('D15').Value = '2.360.88'
$ws.Range('E15').Value = '  +1.10%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '5.42'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -3.67%  '
$ws.Range('D17').Value = '2.043.67'
$ws.Range('E17').Value = '  +0.44%  '
$ws.Range('D18').Value = '36.485.80'
$ws.Range('E18').Value = '  -1.08%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '16.77'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -7.12%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '72.07'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -1.98%  '
$ws.Range('D21').Value = '0.0₃0859'
$ws.Range('E21').Value = '  -2.88%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '238.17'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +1.28%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '5.26'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -1.55%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.36'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -3.45%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '2.13'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +1.90%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '9.29'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -2.92%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '164.72'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -2.49%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '20.12'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +1.32%  '
$ws.Range('E30').Value = '  -0.94%  '
$ws.Range('E31').Value = '  +10.69%  '
$ws.Range('E32').Value = '  -6.16%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '4.49'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -3.94%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.0597'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -2.13%  '
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('E36').Value = '  -0.68%  '
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '2.21'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -0.71%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.0822'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -5.04%  '
$ws.Range('E39').Value = '  -3.25%  '
$ws.Range('E40').Value = '  -5.54%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.0216'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -2.24%  '
$ws.Range('E42').Value = '  -7.50%  '
$ws.Range('E43').Value = '  -2.04%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '94.33'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -2.37%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '1.420.39'
$ws.Range('E45').Value = '  +10.15%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.0912'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -6.06%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '7.63'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +14.49%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '15.92'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -5.21%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '2.86'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +0.61%  '
$ws.Range('E50').Value = '  -2.27%  '
$ws.Range('D51').Value = '2.247.47'
$ws.Range('E51').Value = '  +1.23%  '
